$d = $word.ActiveDocument

$d.Content.Find.Execute("221÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "294÷2=", 2)
$d.Content.Find.Execute("644÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "112÷6=", 2)
$d.Content.Find.Execute("684÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "833÷3=", 2)
$d.Content.Find.Execute("871÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "646÷3=", 2)
$d.Content.Find.Execute("905÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "156÷6=", 2)
$d.Content.Find.Execute("589÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "633÷7=", 2)
$d.Content.Find.Execute("518÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "840÷7=", 2)
$d.Content.Find.Execute("768÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "813÷4=", 2)
$d.Content.Find.Execute("503÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "208÷4=", 2)
$d.Content.Find.Execute("655÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "480÷7=", 2)
$d.Content.Find.Execute("810÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "496÷6=", 2)
$d.Content.Find.Execute("782÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "475÷4=", 2)
$d.Content.Find.Execute("537÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "579÷7=", 2)
$d.Content.Find.Execute("667÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "749÷2=", 2)
$d.Content.Find.Execute("800÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "592÷7=", 2)
$d.Content.Find.Execute("722÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "559÷4=", 2)
$d.Content.Find.Execute("659÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "183÷6=", 2)
$d.Content.Find.Execute("200÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "189÷4=", 2)
$d.Content.Find.Execute("871÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "608÷4=", 2)
$d.Content.Find.Execute("433÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "493÷3=", 2)
$d.Content.Find.Execute("434÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "794÷6=", 2)
$d.Content.Find.Execute("355÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "983÷8=", 2)
$d.Content.Find.Execute("631÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "137÷7=", 2)
$d.Content.Find.Execute("354÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "473÷3=", 2)
$d.Content.Find.Execute("227÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "716÷2=", 2)
